# Auto-generated edit script: refresh cached market-board profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8400.5
$ws.Range("I19").Value = 626.36365
$ws.Range("J19").Value = 14101.533
$ws.Range("K19").Value = 626.36365
$ws.Range("L19").Value = 14101.533
$ws.Range("M19").Value = -451.36365
$ws.Range("N19").Value = -14451.533

$ws.Range("H74").Value = 7125.5
$ws.Range("I74").Value = 5166.6665
$ws.Range("J74").Value = 13002
$ws.Range("K74").Value = 5166.6665
$ws.Range("L74").Value = 13002
$ws.Range("M74").Value = -4230.6665
$ws.Range("N74").Value = -14874

$ws.Range("H77").Value = 7125.5
$ws.Range("I77").Value = 5166.6665
$ws.Range("J77").Value = 13002
$ws.Range("K77").Value = 25833.3325
$ws.Range("L77").Value = 65010
$ws.Range("M77").Value = -21153.3325
$ws.Range("N77").Value = -74370

$ws.Range("H112").Value = 8248.058000000001
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 9439.4
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 28318.2
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -30534.2

$ws.Range("H135").Value = 1728
$ws.Range("I135").Value = 1654.909
$ws.Range("K135").Value = 14894.181
$ws.Range("M135").Value = -12359.181

$ws.Range("H138").Value = 211520.92
$ws.Range("I138").Value = 3337.4167
$ws.Range("J138").Value = 273975.97
$ws.Range("K138").Value = 10012.2501
$ws.Range("L138").Value = 821927.9099999999
$ws.Range("M138").Value = -4872.250100000001
$ws.Range("N138").Value = -832207.9099999999

$ws.Range("H141").Value = 3672.3
$ws.Range("I141").Value = 1838.8462
$ws.Range("J141").Value = 7077.2856
$ws.Range("K141").Value = 5516.5386
$ws.Range("L141").Value = 21231.8568
$ws.Range("M141").Value = -336.5385999999999
$ws.Range("N141").Value = -31591.8568


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 47444.332
$ws.Range("I21").Value = 12407.5
$ws.Range("J21").Value = 64962.75
$ws.Range("K21").Value = 12407.5
$ws.Range("L21").Value = 64962.75
$ws.Range("M21").Value = -12033.5
$ws.Range("N21").Value = -65710.75

$ws.Range("H32").Value = 670449.9399999999
$ws.Range("I32").Value = 765137.4
$ws.Range("J32").Value = 23419
$ws.Range("K32").Value = 765137.4
$ws.Range("L32").Value = 23419
$ws.Range("M32").Value = -764850.4
$ws.Range("N32").Value = -23993

$ws.Range("H63").Value = 6309.3076
$ws.Range("I63").Value = 3691.8
$ws.Range("J63").Value = 7945.25
$ws.Range("K63").Value = 3691.8
$ws.Range("L63").Value = 7945.25
$ws.Range("M63").Value = -3005.8
$ws.Range("N63").Value = -9317.25

$ws.Range("H66").Value = 6309.3076
$ws.Range("I66").Value = 3691.8
$ws.Range("J66").Value = 7945.25
$ws.Range("K66").Value = 18459
$ws.Range("L66").Value = 39726.25
$ws.Range("M66").Value = -15027
$ws.Range("N66").Value = -46590.25

$ws.Range("H131").Value = 39525
$ws.Range("J131").Value = 39525
$ws.Range("L131").Value = 39525
$ws.Range("N131").Value = -49605

$ws.Range("H132").Value = 4582.161
$ws.Range("I132").Value = 4169.3335
$ws.Range("J132").Value = 4969.1875
$ws.Range("K132").Value = 12508.0005
$ws.Range("L132").Value = 14907.5625
$ws.Range("M132").Value = -9978.000499999998
$ws.Range("N132").Value = -19967.5625

$ws.Range("H134").Value = 75429
$ws.Range("J134").Value = 75429
$ws.Range("L134").Value = 75429
$ws.Range("N134").Value = -85569


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1635.1666
$ws.Range("I20").Value = 1711.8572
$ws.Range("J20").Value = 1527.8
$ws.Range("K20").Value = 1711.8572
$ws.Range("L20").Value = 1527.8
$ws.Range("M20").Value = -1464.8572
$ws.Range("N20").Value = -2021.8

$ws.Range("H21").Value = 57500
$ws.Range("J21").Value = 57500
$ws.Range("L21").Value = 57500
$ws.Range("N21").Value = -57972

$ws.Range("H28").Value = 29900
$ws.Range("J28").Value = 29900
$ws.Range("L28").Value = 29900
$ws.Range("N28").Value = -30488

$ws.Range("H62").Value = 90000
$ws.Range("J62").Value = 90000
$ws.Range("L62").Value = 90000
$ws.Range("N62").Value = -91372

$ws.Range("H65").Value = 90000
$ws.Range("J65").Value = 90000
$ws.Range("L65").Value = 270000
$ws.Range("N65").Value = -276864

$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101622

$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -308112

$ws.Range("H92").Value = 68400
$ws.Range("J92").Value = 68400
$ws.Range("L92").Value = 68400
$ws.Range("N92").Value = -73392

$ws.Range("H107").Value = 1479.2222
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 1644.7142
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 1644.7142
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -5484.7142


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 26180
$ws.Range("J10").Value = 31475
$ws.Range("L10").Value = 31475
$ws.Range("N10").Value = -31753

$ws.Range("H16").Value = 873.8
$ws.Range("I16").Value = 780.7
$ws.Range("J16").Value = 1060
$ws.Range("K16").Value = 780.7
$ws.Range("L16").Value = 1060
$ws.Range("M16").Value = -493.7
$ws.Range("N16").Value = -1634

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H31").Value = 4912.1177
$ws.Range("I31").Value = 1276.8889
$ws.Range("J31").Value = 9001.75
$ws.Range("K31").Value = 1276.8889
$ws.Range("L31").Value = 9001.75
$ws.Range("M31").Value = -981.8888999999999
$ws.Range("N31").Value = -9591.75

$ws.Range("H34").Value = 4912.1177
$ws.Range("I34").Value = 1276.8889
$ws.Range("J34").Value = 9001.75
$ws.Range("K34").Value = 1276.8889
$ws.Range("L34").Value = 9001.75
$ws.Range("M34").Value = -1074.8889
$ws.Range("N34").Value = -9405.75

$ws.Range("H113").Value = 873.8
$ws.Range("I113").Value = 780.7
$ws.Range("J113").Value = 1060
$ws.Range("K113").Value = 780.7
$ws.Range("L113").Value = 1060
$ws.Range("M113").Value = 1389.3
$ws.Range("N113").Value = -5400


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1533.0769
$ws.Range("J39").Value = 1888
$ws.Range("L39").Value = 5664
$ws.Range("N39").Value = -6252

$ws.Range("H94").Value = 5075.6665
$ws.Range("J94").Value = 5954.4287
$ws.Range("L94").Value = 17863.2861
$ws.Range("N94").Value = -19215.2861

$ws.Range("H96").Value = 4666.6665
$ws.Range("J96").Value = 4666.6665
$ws.Range("L96").Value = 13999.9995
$ws.Range("N96").Value = -18117.9995

$ws.Range("H126").Value = 5761.4
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5761.4
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 17284.2
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -27164.2

$ws.Range("H134").Value = 5623.724
$ws.Range("I134").Value = 3161.875
$ws.Range("J134").Value = 8653.691999999999
$ws.Range("K134").Value = 9485.625
$ws.Range("L134").Value = 25961.076
$ws.Range("M134").Value = -4415.625
$ws.Range("N134").Value = -36101.076

$ws.Range("H139").Value = 3743.8333
$ws.Range("I139").Value = 1941.4615
$ws.Range("K139").Value = 5824.3845
$ws.Range("M139").Value = -684.3845000000001

$ws.Range("H140").Value = 2163.8235
$ws.Range("I140").Value = 1462.9166
$ws.Range("J140").Value = 3846
$ws.Range("K140").Value = 4388.7498
$ws.Range("L140").Value = 11538
$ws.Range("M140").Value = 791.2502000000004
$ws.Range("N140").Value = -21898

$ws.Range("H141").Value = 4368.1304
$ws.Range("I141").Value = 3804.889
$ws.Range("J141").Value = 6395.8
$ws.Range("K141").Value = 11414.667
$ws.Range("L141").Value = 19187.4
$ws.Range("M141").Value = -6234.667000000001
$ws.Range("N141").Value = -29547.4


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5422.0205
$ws.Range("I70").Value = 5315.517
$ws.Range("J70").Value = 5576.45
$ws.Range("K70").Value = 5315.517
$ws.Range("L70").Value = 5576.45
$ws.Range("M70").Value = -5045.517
$ws.Range("N70").Value = -6116.45

$ws.Range("H73").Value = 5422.0205
$ws.Range("I73").Value = 5315.517
$ws.Range("J73").Value = 5576.45
$ws.Range("K73").Value = 5315.517
$ws.Range("L73").Value = 5576.45
$ws.Range("M73").Value = -4379.517
$ws.Range("N73").Value = -7448.45

$ws.Range("H113").Value = 1511.7858
$ws.Range("I113").Value = 1461.3334
$ws.Range("J113").Value = 1602.6
$ws.Range("K113").Value = 1461.3334
$ws.Range("L113").Value = 1602.6
$ws.Range("M113").Value = 708.6666
$ws.Range("N113").Value = -5942.6

$ws.Range("H132").Value = 2191.6
$ws.Range("I132").Value = 1567.8
$ws.Range("J132").Value = 2815.4
$ws.Range("K132").Value = 4703.4
$ws.Range("L132").Value = 8446.200000000001
$ws.Range("M132").Value = -2173.4
$ws.Range("N132").Value = -13506.2


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 46349.25
$ws.Range("J59").Value = 46349.25
$ws.Range("L59").Value = 46349.25
$ws.Range("N59").Value = -47657.25


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 19442.8
$ws.Range("I56").Value = 3000
$ws.Range("J56").Value = 23553.5
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 23553.5
$ws.Range("M56").Value = -2286
$ws.Range("N56").Value = -24981.5

$ws.Range("H136").Value = 4339.2383
$ws.Range("I136").Value = 4259.2144
$ws.Range("J136").Value = 4499.2856
$ws.Range("K136").Value = 12777.6432
$ws.Range("L136").Value = 13497.8568
$ws.Range("M136").Value = -10227.6432
$ws.Range("N136").Value = -18597.8568

